# Regression Suite for CDS
# Update the TSV/Web data filenames in row 2 to reference the new B Lymphoblastic
# Leukemia test case instead of the old Control test case, and remove the now
# redundant duplicate filename cells from rows 3 and 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "TC14_CDS_phs003164_PrimDiag_B Lymphoblastic Leukemia_TSVData.xlsx"
$ws.Range("E2").Value = "TC14_CDS_phs003164_PrimDiag_B Lymphoblastic Leukemia_WebData.xlsx"

$ws.Range("D3").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("E4").ClearContents()

# Widen column E to fit the new (longer) file name text
$ws.Columns.Item(5).ColumnWidth = 67.3

# Update the view: scroll so row 2 is at the top and select D2
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D2").Select()
